$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear columns F, G, H first so no stale values from the old layout remain,
# then rewrite every populated cell with its final target value. This sheet
# renders a set of markdown tables at varying indentation levels (which
# column A-E holds the "Path" determines nesting depth), followed by
# "Description" (col F) and "More" (col G) columns.
$ws.Range("F1:H23").ClearContents() | Out-Null

$ws.Range("A1").Value = "Path"
$ws.Range("B1").Value = " "
$ws.Range("C1").Value = " "
$ws.Range("D1").Value = " "
$ws.Range("E1").Value = " "
$ws.Range("F1").Value = "Introduction"
$ws.Range("G1").Value = "More"
$ws.Range("A2").Value = "[/config](/config)"
$ws.Range("F2").Value = "Search tag configuration.  Partial ML Gradle configurations, referenced during custom token replacement."
$ws.Range("G2").Value = "[Custom Token Replacement](/docs/lux-backend-deployment.md#custom-token-replacement)"
$ws.Range("A3").Value = "[/docs](/docs)"
$ws.Range("F3").Value = "Majority of the backend's documentation."
$ws.Range("G3").Value = "[/README.md](/README.md)"
$ws.Range("A4").Value = "[/postman](/postman)"
$ws.Range("F4").Value = "Contains exports of Postman LUX-related requests and environment template."
$ws.Range("G4").Value = "[LUX Postman Workspace](/docs/lux-postman-workspace.md)"
$ws.Range("A5").Value = "[/scripts](/scripts)"
$ws.Range("F5").Value = "Developer and admin scripts not deployed to an environment.  May be executed from within VS Code, and thus can serve as a way to collaborate on queries outside of a Query Console workspace."
$ws.Range("A6").Value = "[/src](/src)"
$ws.Range("B7").Value = "[/main](/src/main)"
$ws.Range("F7").Value = "All of the project's runtime code, most of its MarkLogic configuration, and some data all within ML Gradle's conventions."
$ws.Range("C8").Value = "[/ml-config](/src/main/ml-config)"
$ws.Range("F8").Value = "All of the project's ML Gradle configuration directories.  Selected ones may vary by environment."
$ws.Range("G8").Value = "[Gradle Properties](/docs/lux-backend-deployment.md#gradle-properties)"
$ws.Range("D9").Value = "[/base](/src/main/ml-config/base)"
$ws.Range("F9").Value = "The base configuration directory applicable to all environments.  It includes the group configuration, main content database, roles, and application servers."
$ws.Range("D10").Value = "[/base-secured](/src/main/ml-config/base-secured)"
$ws.Range("F10").Value = "HTTPS settings that stack on top of the base configuration."
$ws.Range("D11").Value = "[/base-unsecured](/src/main/ml-config/base-unsecured)"
$ws.Range("F11").Value = "Defines a local, non-admin user to perform most of deployments with."
$ws.Range("C12").Value = "[/ml-data](/src/main/ml-data)"
$ws.Range("F12").Value = "ML Gradle's default data directory. Presently only used for thesauri but need not be limited to. Not expecting to load datasets here though."
$ws.Range("C13").Value = "[/ml-modules](/src/main/ml-modules)"
$ws.Range("F13").Value = "The modules applicable to all environments."
$ws.Range("D14").Value = "[/options](/src/main/ml-modules/options)"
$ws.Range("F14").Value = "/v1/search options, which this project does not use."
$ws.Range("D15").Value = "[/root](/src/main/ml-modules/root)"
$ws.Range("E16").Value = "[/config](/src/main/ml-modules/root/config)"
$ws.Range("F16").Value = "Configuration for search, facets, and more.  Includes placeholder files that are replaced during deployment."
$ws.Range("E17").Value = "[/data](/src/main/ml-modules/root/data)"
$ws.Range("F17").Value = "Includes the words to exclude from search criteria."
$ws.Range("E18").Value = "[/ds](/src/main/ml-modules/root/ds)"
$ws.Range("F18").Value = "All of LUX's custom MarkLogic data services, which should just be wrappers to library modules."
$ws.Range("G18").Value = "[LUX Backend API Usage Documentation](/docs/lux-backend-api-usage.md)"
$ws.Range("E19").Value = "[/lib](/src/main/ml-modules/root/lib)"
$ws.Range("F19").Value = "The heart of LUX's backend implementation, where developers get to spend most of their time when they're lucky :)"
$ws.Range("E20").Value = "[/runDuringDeployment](/src/main/ml-modules/root/runDuringDeployment)"
$ws.Range("F20").Value = "Includes scripts to deploy then execute during deployment, directly supporting generators for the remaining search terms, related lists, and advanced search configuration."
$ws.Range("G20").Value = "[LUX Gradle Tasks](/docs/lux-backend-build-tool-and-tasks.md#lux-gradle-tasks)"
$ws.Range("E21").Value = "[/utils](/src/main/ml-modules/root/utils)"
$ws.Range("F21").Value = "A few utility functions and classes used by the library code."
$ws.Range("C22").Value = "[/templates](/src/main/templates)"
$ws.Range("F22").Value = "JavaScript template files used by [/build.gradle](/build.gradle)"
$ws.Range("G22").Value = "[JavaScript Template Files](/docs/lux-backend-build-tool-and-tasks.md#javascript-template-files)"
$ws.Range("A23").Value = "[/build.gradle](/build.gradle)"
$ws.Range("F23").Value = "The build script."
$ws.Range("G23").Value = "[LUX Backend Local Developer Environment](/docs/lux-backend-setup-local-env.md), [LUX Backend Deployment](/docs/lux-backend-deployment.md), [LUX Backend Build Tool and Tasks](/docs/lux-backend-build-tool-and-tasks.md)"
